# Project Allocation - add Week2 section and rename a few tasks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert the rows we need first (shift everything below down) ---
# Row 4: "Week1" divider row (new)
$ws.Rows("4:4").Insert()
# Row 9: "Week2" divider row (new)
$ws.Rows("9:9").Insert()
# Rows 10-13: Week2 data row + 3 blank placeholder rows (new)
$ws.Rows("10:13").Insert()

# --- 2. Update the Week1 task text that changed ---
$ws.Range("B5").Value = "Invoices user front + back"
$ws.Range("B7").Value = "Dashboard administrator front + back"

# --- 3. Week1 header row (row 4) ---
$ws.Range("A4").Value = "Week1"
$ws.Range("A4").Style = $ws.Range("A3").Style
$ws.Range("B4").Style = $ws.Range("B3").Style
$ws.Range("C4").Style = $ws.Range("C3").Style

# --- 4. Week2 header row (row 9), bold black font ---
$ws.Range("A9").Value = "Week2"
$ws.Range("A9").Font.Bold = $true
$ws.Range("A9").Font.Color = 0

# --- 5. Week2 data row (row 10) ---
$ws.Range("A10").Value = "Marcu Andrei"
$ws.Range("B10").Value = "Login front + back, routing and login protection"
$ws.Range("C10").Value = "7 zile"

# --- 6. Blank placeholder rows under Week2 (rows 11-13), bold black font like row 9 ---
$ws.Range("A11").Font.Bold = $true
$ws.Range("A11").Font.Color = 0
$ws.Range("A12").Font.Bold = $true
$ws.Range("A12").Font.Color = 0
$ws.Range("A13").Font.Bold = $true
$ws.Range("A13").Font.Color = 0

# --- 7. Sheet view / selection ---
$ws.Range("B13").Select()

# --- 8. Page setup (portrait) ---
$ws.PageSetup.Orientation = 1
